# Update the "Förändrad" date column (C) from 2026-02-22 (46075) to
# 2026-02-23 (46076) for every data row (rows 2 through 121) on the
# "Avverkningsanmälningar" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 121 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value = 46076
    }
}
